# Update the date heading (first paragraph)
$d = $word.ActiveDocument
$d.Paragraphs.Item(1).Range.Text = "2024-04-11 Thursday"

# Update each arithmetic-problem table cell (20 rows x 5 cols), in document order
$tbl = $d.Tables.Item(1)

$tbl.Cell(1, 1).Range.Text = "3+81="
$tbl.Cell(1, 2).Range.Text = "88+3="
$tbl.Cell(1, 3).Range.Text = "56+4="
$tbl.Cell(1, 4).Range.Text = "8+13="
$tbl.Cell(1, 5).Range.Text = "57-15="
$tbl.Cell(2, 1).Range.Text = "69-2="
$tbl.Cell(2, 2).Range.Text = "33+64="
$tbl.Cell(2, 3).Range.Text = "45-45="
$tbl.Cell(2, 4).Range.Text = "3+81="
$tbl.Cell(2, 5).Range.Text = "53-17="
$tbl.Cell(3, 1).Range.Text = "99-59="
$tbl.Cell(3, 2).Range.Text = "41+33="
$tbl.Cell(3, 3).Range.Text = "38-7="
$tbl.Cell(3, 4).Range.Text = "40+14="
$tbl.Cell(3, 5).Range.Text = "43-8="
$tbl.Cell(4, 1).Range.Text = "52+23="
$tbl.Cell(4, 2).Range.Text = "12-2="
$tbl.Cell(4, 3).Range.Text = "66+20="
$tbl.Cell(4, 4).Range.Text = "90-75="
$tbl.Cell(4, 5).Range.Text = "15+64="
$tbl.Cell(5, 1).Range.Text = "45-29="
$tbl.Cell(5, 2).Range.Text = "54+6="
$tbl.Cell(5, 3).Range.Text = "40+9="
$tbl.Cell(5, 4).Range.Text = "0+91="
$tbl.Cell(5, 5).Range.Text = "2+2="
$tbl.Cell(6, 1).Range.Text = "75+13="
$tbl.Cell(6, 2).Range.Text = "29+3="
$tbl.Cell(6, 3).Range.Text = "91-89="
$tbl.Cell(6, 4).Range.Text = "40+9="
$tbl.Cell(6, 5).Range.Text = "99-86="
$tbl.Cell(7, 1).Range.Text = "34+63="
$tbl.Cell(7, 2).Range.Text = "61+5="
$tbl.Cell(7, 3).Range.Text = "53+5="
$tbl.Cell(7, 4).Range.Text = "33-11="
$tbl.Cell(7, 5).Range.Text = "91-79="
$tbl.Cell(8, 1).Range.Text = "39-9="
$tbl.Cell(8, 2).Range.Text = "66-58="
$tbl.Cell(8, 3).Range.Text = "79+8="
$tbl.Cell(8, 4).Range.Text = "29-10="
$tbl.Cell(8, 5).Range.Text = "28+50="
$tbl.Cell(9, 1).Range.Text = "18+15="
$tbl.Cell(9, 2).Range.Text = "34+44="
$tbl.Cell(9, 3).Range.Text = "41-38="
$tbl.Cell(9, 4).Range.Text = "92-2="
$tbl.Cell(9, 5).Range.Text = "88-36="
$tbl.Cell(10, 1).Range.Text = "53-13="
$tbl.Cell(10, 2).Range.Text = "75-46="
$tbl.Cell(10, 3).Range.Text = "84+13="
$tbl.Cell(10, 4).Range.Text = "57+22="
$tbl.Cell(10, 5).Range.Text = "46-26="
$tbl.Cell(11, 1).Range.Text = "46+6="
$tbl.Cell(11, 2).Range.Text = "51+9="
$tbl.Cell(11, 3).Range.Text = "65-18="
$tbl.Cell(11, 4).Range.Text = "84+5="
$tbl.Cell(11, 5).Range.Text = "87-83="
$tbl.Cell(12, 1).Range.Text = "31+27="
$tbl.Cell(12, 2).Range.Text = "13+16="
$tbl.Cell(12, 3).Range.Text = "21+62="
$tbl.Cell(12, 4).Range.Text = "96-78="
$tbl.Cell(12, 5).Range.Text = "4+91="
$tbl.Cell(13, 1).Range.Text = "89-85="
$tbl.Cell(13, 2).Range.Text = "37+44="
$tbl.Cell(13, 3).Range.Text = "74+4="
$tbl.Cell(13, 4).Range.Text = "0+49="
$tbl.Cell(13, 5).Range.Text = "0+84="
$tbl.Cell(14, 1).Range.Text = "70-13="
$tbl.Cell(14, 2).Range.Text = "0+11="
$tbl.Cell(14, 3).Range.Text = "25+43="
$tbl.Cell(14, 4).Range.Text = "92-89="
$tbl.Cell(14, 5).Range.Text = "32-1="
$tbl.Cell(15, 1).Range.Text = "83-7="
$tbl.Cell(15, 2).Range.Text = "47-11="
$tbl.Cell(15, 3).Range.Text = "68-59="
$tbl.Cell(15, 4).Range.Text = "8+34="
$tbl.Cell(15, 5).Range.Text = "58+5="
$tbl.Cell(16, 1).Range.Text = "86-47="
$tbl.Cell(16, 2).Range.Text = "32+26="
$tbl.Cell(16, 3).Range.Text = "72-71="
$tbl.Cell(16, 4).Range.Text = "91-7="
$tbl.Cell(16, 5).Range.Text = "79-42="
$tbl.Cell(17, 1).Range.Text = "78-34="
$tbl.Cell(17, 2).Range.Text = "8+52="
$tbl.Cell(17, 3).Range.Text = "91-75="
$tbl.Cell(17, 4).Range.Text = "97-50="
$tbl.Cell(17, 5).Range.Text = "71+16="
$tbl.Cell(18, 1).Range.Text = "30+11="
$tbl.Cell(18, 2).Range.Text = "18+5="
$tbl.Cell(18, 3).Range.Text = "3+67="
$tbl.Cell(18, 4).Range.Text = "19+47="
$tbl.Cell(18, 5).Range.Text = "7+49="
$tbl.Cell(19, 1).Range.Text = "25+9="
$tbl.Cell(19, 2).Range.Text = "47-27="
$tbl.Cell(19, 3).Range.Text = "33-9="
$tbl.Cell(19, 4).Range.Text = "33-8="
$tbl.Cell(19, 5).Range.Text = "87-17="
$tbl.Cell(20, 1).Range.Text = "66+32="
$tbl.Cell(20, 2).Range.Text = "34+22="
$tbl.Cell(20, 3).Range.Text = "43+53="
$tbl.Cell(20, 4).Range.Text = "95-4="
$tbl.Cell(20, 5).Range.Text = "22+76="
